$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Heading: "3" + "/7/19 ELEN4012 Planning meeting" -> single run
#    A Find/Replace that spans both runs causes Word to coalesce them
#    into one run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("3/7/19 ELEN4012 Planning meeting", $true, $false, $false, $false, $false, $true, 1, $false, "3/7/19 ELEN4012 Planning meeting", 2) | Out-Null

# ------------------------------------------------------------------
# 2. "ML exploration time: best " + bookmark _GoBack + "tasks? Software?"
#    -> single run, bookmark removed.
# ------------------------------------------------------------------
$d.Content.Find.Execute("ML exploration time: best tasks? Software?", $true, $false, $false, $false, $false, $true, 1, $false, "ML exploration time: best tasks? Software?", 2) | Out-Null

# ------------------------------------------------------------------
# 3. "Too many options!" -> "Tutorials to get familiar" (same paragraph)
# ------------------------------------------------------------------
$d.Content.Find.Execute("Too many options!", $true, $false, $false, $false, $false, $true, 1, $false, "Tutorials to get familiar", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Insert the four new bullet paragraphs right after "Tutorials to
#    get familiar", reproducing the exact run/proofErr/bookmark layout.
#    The document currently ends in a single trailing blank <w:p/>
#    right before the sectPr; InsertXML on a range collapsed to the
#    start of a paragraph *replaces* that paragraph, so we replace the
#    trailing blank paragraph with [4 new paragraphs + a fresh blank
#    paragraph] in one shot, which keeps the original trailing blank
#    paragraph in place logically.
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Tutorials to get familiar`r") {
        $targetPara = $p
    }
}
$trailingBlank = $targetPara.Next()
$insertRange = $trailingBlank.Range
$insertRange.Collapse(1)

$newParagraphsXml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="ListBullet"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r><w:t>Discussion of image processing:</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="ListBullet"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Sorbel</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>/Edge detection</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="ListBullet"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r><w:t>Motion vectors?</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="ListBullet"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Pycine</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> : openCV</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
"@

$insertRange.InsertXML($newParagraphsXml)

# ------------------------------------------------------------------
# 5. Mint footnotes.xml / endnotes.xml parts (with just the default
#    separator / continuationSeparator boilerplate) the same way Word
#    does the first time a document's footnote machinery is touched.
# ------------------------------------------------------------------
$tempNote = $d.Footnotes.Add($d.Paragraphs(1).Range, "", "x")
$d.Footnotes(1).Delete()
